# Update the stringer position data on "Cross Section Properties" sheet.
# Rows 6 & 7 (Top Stringer1 / Top Stringer2) get real z'i (E) values instead
# of the placeholder "TBD" text, which lets the downstream K/L/M/N/O/...
# formulas (previously #VALUE!) compute real numbers. Row 7's y'i (D) also
# switches from the old shared "(1/16)+(1/8)/2" formula to a relative
# reference off row 6, and rows 11-14 (Bottom Stringer1..4) get updated
# y'i / z'i coordinates. The TOTALS row (16) sums are widened to include
# the now-valid rows 6 & 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cross Section Properties")

# Row 6 - Top Stringer1: E6 becomes a plain numeric value (was text "TBD")
$ws.Range("E6").Value = -0.0625

# Row 7 - Top Stringer2
$ws.Range("D7").Formula = "=D6+0.125"
$ws.Range("E7").Value = -0.0625

# Row 11 - Bottom Stringer1: now formulas, and loses the "s=5" fill style
$ws.Range("D11").Formula = "=-1.6046754518"
$ws.Range("E11").Formula = "=-0.092597515"
$ws.Range("D11:E11").Interior.ColorIndex = 0
$ws.Range("D11:E11").Style = "Normal"

# Row 12 - Bottom Stringer2: now formulas, and loses the "s=5" fill style
$ws.Range("D12").Formula = "=-1.7254548192"
$ws.Range("E12").Formula = "=-0.1248053464"
$ws.Range("D12:E12").Interior.ColorIndex = 0
$ws.Range("D12:E12").Style = "Normal"

# Row 13 - Bottom Stringer3 (keeps its existing fill style)
$ws.Range("D13").Formula = "=-1.8462341867"
$ws.Range("E13").Formula = "=-0.1570131777"

# Row 14 - Bottom Stringer4 (keeps its existing fill style)
$ws.Range("D14").Formula = "=-1.9670135541"
$ws.Range("E14").Formula = "=-0.189221009"

# TOTALS row (16): widen the summed ranges to include the now-valid rows 6 & 7
$ws.Range("C16").Formula = "=SUM(C2:C7,C11:C14)"
$ws.Range("F16").Formula = "=SUM(F2:F7,F11:F14)"
$ws.Range("G16").Formula = "=SUM(G2:G7,G11:G14)"
$ws.Range("H16").Formula = "=SUM(H2:H7,H11:H14)"
$ws.Range("I16").Formula = "=SUM(I2:I7,I11:I14)"
$ws.Range("J16").Formula = "=SUM(J2:J7,J11:J14)"
$ws.Range("M16").Formula = "=SUM(M2:M7,M11:M14)"
$ws.Range("N16").Formula = "=SUM(N2:N7,N11:N14)"
$ws.Range("O16").Formula = "=SUM(O2:O7,O11:O14)"
$ws.Range("R16").Formula = "=SUM(R2:R7,R11:R14)"
$ws.Range("S16").Formula = "=SUM(S2:S7,S11:S14)"
$ws.Range("T16").Formula = "=SUM(T2:T7,T11:T14)"
$ws.Range("W16").Formula = "=SUM(W2:W7,W11:W14)"
$ws.Range("X16").Formula = "=SUM(X2:X7,X11:X14)"
$ws.Range("Y16").Formula = "=SUM(Y2:Y7,Y11:Y14)"
$ws.Range("AZ16").Formula = "=SUM(AZ2:AZ7,AZ11:AZ14)"
$ws.Range("BA16").Formula = "=SUM(BA2:BA7,BA11:BA14)"
$ws.Range("BB16").Formula = "=SUM(BB2:BB7,BB11:BB14)"

$excel.CalculateFullRebuild()

# Work around stale cached shared-formula results for the cells that moved
# out of an #VALUE! state (rows 6 & 7) by re-entering their formulas so the
# engine re-evaluates them fresh.
$ws.Range("L6").Formula = '=E6-$AI$3'
$ws.Range("L7").Formula = '=E7-$AI$3'
$ws.Range("V6").Formula = '=E6-$AI$5'
$ws.Range("V7").Formula = '=E7-$AI$5'
$ws.Range("N6").Formula = "=C6*L6^2"
$ws.Range("N7").Formula = "=C7*L7^2"
$ws.Range("O6").Formula = "=C6*K6*L6"
$ws.Range("O7").Formula = "=C7*K7*L7"
$ws.Range("X6").Formula = "=C6*V6^2"
$ws.Range("X7").Formula = "=C7*V7^2"
$ws.Range("Y6").Formula = "=C6*U6*V6"
$ws.Range("Y7").Formula = "=C7*U7*V7"
$ws.Range("AZ6").Formula = "=(Q6)*(H6+X6)"
$ws.Range("AZ7").Formula = "=(Q7)*(H7+X7)"
$ws.Range("BA6").Formula = "=(Q6)*(I6+W6)"
$ws.Range("BA7").Formula = "=(Q7)*(I7+W7)"
$ws.Range("BB6").Formula = "=(Q6)*(J6+Y6)"
$ws.Range("BB7").Formula = "=(Q7)*(J7+Y7)"

$excel.CalculateFullRebuild()
